$wb = $excel.ActiveWorkbook

# Updated Leve profit/price figures (columns H-N) produced by the scheduled
# market-data runner, applied per character sheet below.

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1297.6666
$ws.Range("I106").Value = 946.5
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 946.5
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -315.5
$ws.Range("N106").Value = -3262
$ws.Range("H115").Value = 1933.75
$ws.Range("I115").Value = 1051.25
$ws.Range("J115").Value = 2375
$ws.Range("K115").Value = 3153.75
$ws.Range("L115").Value = 7125
$ws.Range("M115").Value = -1586.75
$ws.Range("N115").Value = -10259
$ws.Range("H129").Value = 4808698.5
$ws.Range("I129").Value = 41667964
$ws.Range("J129").Value = 968.36957
$ws.Range("K129").Value = 125003892
$ws.Range("L129").Value = 2905.10871
$ws.Range("M129").Value = -124998892
$ws.Range("N129").Value = -12905.10871

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1325.26
$ws.Range("I32").Value = 1191.2556
$ws.Range("J32").Value = 2531.3
$ws.Range("K32").Value = 1191.2556
$ws.Range("L32").Value = 2531.3
$ws.Range("M32").Value = -904.2556
$ws.Range("N32").Value = -3105.3
$ws.Range("H61").Value = 3032.3547
$ws.Range("I61").Value = 1026.8667
$ws.Range("K61").Value = 1026.8667
$ws.Range("M61").Value = -814.8667
$ws.Range("H74").Value = 823.13336
$ws.Range("I74").Value = 765.1539
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 765.1539
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = 108.8461
$ws.Range("N74").Value = -2948
$ws.Range("H77").Value = 823.13336
$ws.Range("I77").Value = 765.1539
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 3825.7695
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = 542.2304999999997
$ws.Range("N77").Value = -14736
$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140
$ws.Range("H136").Value = 3032.3547
$ws.Range("I136").Value = 1026.8667
$ws.Range("K136").Value = 3080.6001
$ws.Range("M136").Value = -530.6001000000001
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 8749.25
$ws.Range("I19").Value = 7499
$ws.Range("J19").Value = 9999.5
$ws.Range("K19").Value = 7499
$ws.Range("L19").Value = 9999.5
$ws.Range("M19").Value = -7326
$ws.Range("N19").Value = -10345.5
$ws.Range("H107").Value = 3967.182
$ws.Range("J107").Value = 5505
$ws.Range("L107").Value = 5505
$ws.Range("N107").Value = -9345
$ws.Range("H134").Value = 4141.8335
$ws.Range("I134").Value = 3027.75
$ws.Range("K134").Value = 9083.25
$ws.Range("M134").Value = -6548.25
$ws.Range("H141").Value = 35426.43
$ws.Range("J141").Value = 29597
$ws.Range("L141").Value = 29597
$ws.Range("N141").Value = -39957

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2668.52
$ws.Range("I31").Value = 1485.125
$ws.Range("J31").Value = 4772.3335
$ws.Range("K31").Value = 1485.125
$ws.Range("L31").Value = 4772.3335
$ws.Range("M31").Value = -1190.125
$ws.Range("N31").Value = -5362.3335
$ws.Range("H34").Value = 2668.52
$ws.Range("I34").Value = 1485.125
$ws.Range("J34").Value = 4772.3335
$ws.Range("K34").Value = 1485.125
$ws.Range("L34").Value = 4772.3335
$ws.Range("M34").Value = -1283.125
$ws.Range("N34").Value = -5176.3335
$ws.Range("H51").Value = 10000
$ws.Range("J51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -11472
$ws.Range("H59").Value = 16931.75
$ws.Range("J59").Value = 16931.75
$ws.Range("L59").Value = 16931.75
$ws.Range("N59").Value = -19221.75
$ws.Range("H61").Value = 10000
$ws.Range("J61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("N61").Value = -10696
$ws.Range("H127").Value = 32993.332
$ws.Range("J127").Value = 32993.332
$ws.Range("L127").Value = 32993.332
$ws.Range("N127").Value = -42913.332
$ws.Range("H141").Value = 29372.727
$ws.Range("J141").Value = 29372.727
$ws.Range("L141").Value = 29372.727
$ws.Range("N141").Value = -39732.727

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1081.2667
$ws.Range("I44").Value = 494.2
$ws.Range("J44").Value = 1374.8
$ws.Range("K44").Value = 1482.6
$ws.Range("L44").Value = 4124.4
$ws.Range("M44").Value = -1084.6
$ws.Range("N44").Value = -4920.4
$ws.Range("H55").Value = 2974.1667
$ws.Range("J55").Value = 3474
$ws.Range("L55").Value = 10422
$ws.Range("N55").Value = -10776
$ws.Range("H122").Value = 1524.25
$ws.Range("J122").Value = 2407.9092
$ws.Range("L122").Value = 21671.1828
$ws.Range("N122").Value = -26571.1828
$ws.Range("H131").Value = 1440
$ws.Range("J131").Value = 1283.125
$ws.Range("L131").Value = 3849.375
$ws.Range("N131").Value = -13929.375
$ws.Range("H133").Value = 4125.294
$ws.Range("I133").Value = 5892.857
$ws.Range("J133").Value = 2888
$ws.Range("K133").Value = 17678.571
$ws.Range("L133").Value = 8664
$ws.Range("M133").Value = -12618.571
$ws.Range("N133").Value = -18784
$ws.Range("H138").Value = 2103.4167
$ws.Range("J138").Value = 2895
$ws.Range("L138").Value = 8685
$ws.Range("N138").Value = -18965
$ws.Range("H141").Value = 3227.2727
$ws.Range("I141").Value = 2562.5
$ws.Range("K141").Value = 7687.5
$ws.Range("M141").Value = -2507.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()
$ws.Range("H132").Value = 3861.3794
$ws.Range("I132").Value = 4001.6924
$ws.Range("J132").Value = 3747.375
$ws.Range("K132").Value = 12005.0772
$ws.Range("L132").Value = 11242.125
$ws.Range("M132").Value = -9475.0772
$ws.Range("N132").Value = -16302.125
$ws.Range("H135").Value = 29333.334
$ws.Range("J135").Value = 29333.334
$ws.Range("L135").Value = 29333.334
$ws.Range("N135").Value = -39473.334
$ws.Range("H137").Value = 29546
$ws.Range("J137").Value = 29546
$ws.Range("L137").Value = 29546
$ws.Range("N137").Value = -39746
$ws.Range("H138").Value = 40437.5
$ws.Range("J138").Value = 40437.5
$ws.Range("L138").Value = 40437.5
$ws.Range("N138").Value = -50717.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3000
$ws.Range("J12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("N12").Value = -3340
$ws.Range("H109").Value = 30600
$ws.Range("J109").Value = 30600
$ws.Range("L109").Value = 30600
$ws.Range("N109").Value = -33374
$ws.Range("H132").Value = 3240.4
$ws.Range("I132").Value = 2069.7693
$ws.Range("J132").Value = 4508.5835
$ws.Range("K132").Value = 6209.3079
$ws.Range("L132").Value = 13525.7505
$ws.Range("M132").Value = -3679.3079
$ws.Range("N132").Value = -18585.7505
$ws.Range("H140").Value = 29551.666
$ws.Range("J140").Value = 29551.666
$ws.Range("L140").Value = 29551.666
$ws.Range("N140").Value = -39911.666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 35839.8
$ws.Range("J135").Value = 35839.8
$ws.Range("L135").Value = 35839.8
$ws.Range("N135").Value = -45979.8
$ws.Range("H138").Value = 26666.666
$ws.Range("J138").Value = 26666.666
$ws.Range("L138").Value = 26666.666
$ws.Range("N138").Value = -36946.666

